$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 3416.9167
$ws.Cells.Item(51, 9).Value = 2112.75
$ws.Cells.Item(51, 10).Value = 4069
$ws.Cells.Item(51, 11).Value = 2112.75
$ws.Cells.Item(51, 12).Value = 4069
$ws.Cells.Item(51, 13).Value = -1628.75
$ws.Cells.Item(51, 14).Value = -5037

$ws.Cells.Item(58, 8).Value = 758.4211
$ws.Cells.Item(58, 9).Value = 444.84616
$ws.Cells.Item(58, 10).Value = 1437.8334
$ws.Cells.Item(58, 11).Value = 1334.53848
$ws.Cells.Item(58, 12).Value = 4313.5002
$ws.Cells.Item(58, 13).Value = -1184.53848
$ws.Cells.Item(58, 14).Value = -4613.5002

$ws.Cells.Item(62, 8).Value = 3389.5
$ws.Cells.Item(62, 9).Value = 2220
$ws.Cells.Item(62, 10).Value = 3839.3076
$ws.Cells.Item(62, 11).Value = 2220
$ws.Cells.Item(62, 12).Value = 3839.3076
$ws.Cells.Item(62, 13).Value = -1596
$ws.Cells.Item(62, 14).Value = -5087.3076

$ws.Cells.Item(65, 8).Value = 3389.5
$ws.Cells.Item(65, 9).Value = 2220
$ws.Cells.Item(65, 10).Value = 3839.3076
$ws.Cells.Item(65, 11).Value = 11100
$ws.Cells.Item(65, 12).Value = 19196.538
$ws.Cells.Item(65, 13).Value = -7980
$ws.Cells.Item(65, 14).Value = -25436.538

$ws.Cells.Item(70, 8).Value = 778.7778
$ws.Cells.Item(70, 9).Value = 600
$ws.Cells.Item(70, 10).Value = 829.8570999999999
$ws.Cells.Item(70, 11).Value = 1800
$ws.Cells.Item(70, 12).Value = 2489.5713
$ws.Cells.Item(70, 13).Value = -1530
$ws.Cells.Item(70, 14).Value = -3029.5713

$ws.Cells.Item(73, 8).Value = 778.7778
$ws.Cells.Item(73, 9).Value = 600
$ws.Cells.Item(73, 10).Value = 829.8570999999999
$ws.Cells.Item(73, 11).Value = 1800
$ws.Cells.Item(73, 12).Value = 2489.5713
$ws.Cells.Item(73, 13).Value = -864
$ws.Cells.Item(73, 14).Value = -4361.5713

$ws.Cells.Item(98, 8).Value = 2035.8667
$ws.Cells.Item(98, 9).Value = 1245.7142
$ws.Cells.Item(98, 10).Value = 2727.25
$ws.Cells.Item(98, 11).Value = 1245.7142
$ws.Cells.Item(98, 12).Value = 2727.25
$ws.Cells.Item(98, 13).Value = 252.2858000000001
$ws.Cells.Item(98, 14).Value = -5723.25

$ws.Cells.Item(122, 8).Value = 2035.8667
$ws.Cells.Item(122, 9).Value = 1245.7142
$ws.Cells.Item(122, 10).Value = 2727.25
$ws.Cells.Item(122, 11).Value = 3737.1426
$ws.Cells.Item(122, 12).Value = 8181.75
$ws.Cells.Item(122, 13).Value = -1287.1426
$ws.Cells.Item(122, 14).Value = -13081.75

$ws.Cells.Item(129, 8).Value = 789.64105
$ws.Cells.Item(129, 9).Value = 298.375
$ws.Cells.Item(129, 10).Value = 916.4194
$ws.Cells.Item(129, 11).Value = 895.125
$ws.Cells.Item(129, 12).Value = 2749.2582
$ws.Cells.Item(129, 13).Value = 4104.875
$ws.Cells.Item(129, 14).Value = -12749.2582

$ws.Cells.Item(132, 8).Value = 8336269
$ws.Cells.Item(132, 9).Value = 11366920
$ws.Cells.Item(132, 10).Value = 1978
$ws.Cells.Item(132, 11).Value = 34100760
$ws.Cells.Item(132, 12).Value = 5934
$ws.Cells.Item(132, 13).Value = -34098230
$ws.Cells.Item(132, 14).Value = -10994

$ws.Cells.Item(138, 8).Value = 4070.4634
$ws.Cells.Item(138, 9).Value = 2143.1143
$ws.Cells.Item(138, 10).Value = 5505.7236
$ws.Cells.Item(138, 11).Value = 6429.342900000001
$ws.Cells.Item(138, 12).Value = 16517.1708
$ws.Cells.Item(138, 13).Value = -1289.342900000001
$ws.Cells.Item(138, 14).Value = -26797.1708

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1038.0278
$ws.Cells.Item(45, 9).Value = 1015.129
$ws.Cells.Item(45, 10).Value = 1180
$ws.Cells.Item(45, 11).Value = 1015.129
$ws.Cells.Item(45, 12).Value = 1180
$ws.Cells.Item(45, 13).Value = -638.129
$ws.Cells.Item(45, 14).Value = -1934

$ws.Cells.Item(97, 8).Value = 1925.6
$ws.Cells.Item(97, 9).Value = 1482.72
$ws.Cells.Item(97, 10).Value = 4140
$ws.Cells.Item(97, 11).Value = 1482.72
$ws.Cells.Item(97, 12).Value = 4140
$ws.Cells.Item(97, 13).Value = -986.72
$ws.Cells.Item(97, 14).Value = -5132

$ws.Cells.Item(109, 8).Value = 22634.25
$ws.Cells.Item(109, 10).Value = 22634.25
$ws.Cells.Item(109, 12).Value = 22634.25
$ws.Cells.Item(109, 14).Value = -25408.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1085.4615
$ws.Cells.Item(134, 9).Value = 995.1081
$ws.Cells.Item(134, 10).Value = 2757
$ws.Cells.Item(134, 11).Value = 2985.3243
$ws.Cells.Item(134, 12).Value = 8271
$ws.Cells.Item(134, 13).Value = -450.3243000000002
$ws.Cells.Item(134, 14).Value = -13341

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 50000
$ws.Cells.Item(37, 10).Value = 50000
$ws.Cells.Item(37, 12).Value = 150000
$ws.Cells.Item(37, 14).Value = -150224

$ws.Cells.Item(122, 8).Value = 1490.625
$ws.Cells.Item(122, 9).Value = 400.5
$ws.Cells.Item(122, 10).Value = 1854
$ws.Cells.Item(122, 11).Value = 3604.5
$ws.Cells.Item(122, 12).Value = 16686
$ws.Cells.Item(122, 13).Value = -1154.5
$ws.Cells.Item(122, 14).Value = -21586

$ws.Cells.Item(139, 8).Value = 1307.7222
$ws.Cells.Item(139, 9).Value = 852.2308
$ws.Cells.Item(139, 10).Value = 2492
$ws.Cells.Item(139, 11).Value = 2556.6924
$ws.Cells.Item(139, 12).Value = 7476
$ws.Cells.Item(139, 13).Value = 2583.3076
$ws.Cells.Item(139, 14).Value = -17756

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2217.2896
$ws.Cells.Item(122, 9).Value = 1113.64
$ws.Cells.Item(122, 10).Value = 4339.6924
$ws.Cells.Item(122, 11).Value = 3340.92
$ws.Cells.Item(122, 12).Value = 13019.0772
$ws.Cells.Item(122, 13).Value = -890.9200000000001
$ws.Cells.Item(122, 14).Value = -17919.0772

$ws.Cells.Item(132, 8).Value = 2401.35
$ws.Cells.Item(132, 9).Value = 1770.3846
$ws.Cells.Item(132, 11).Value = 5311.1538
$ws.Cells.Item(132, 13).Value = -2781.1538

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 435.57693
$ws.Cells.Item(55, 9).Value = 397.11765
$ws.Cells.Item(55, 10).Value = 508.22223
$ws.Cells.Item(55, 11).Value = 397.11765
$ws.Cells.Item(55, 12).Value = 508.22223
$ws.Cells.Item(55, 13).Value = -224.11765
$ws.Cells.Item(55, 14).Value = -854.2222300000001

$ws.Cells.Item(68, 8).Value = 2641.111
$ws.Cells.Item(68, 9).Value = 1142
$ws.Cells.Item(68, 10).Value = 2871.7437
$ws.Cells.Item(68, 11).Value = 1142
$ws.Cells.Item(68, 12).Value = 2871.7437
$ws.Cells.Item(68, 13).Value = -393
$ws.Cells.Item(68, 14).Value = -4369.7437

$ws.Cells.Item(71, 8).Value = 2641.111
$ws.Cells.Item(71, 9).Value = 1142
$ws.Cells.Item(71, 10).Value = 2871.7437
$ws.Cells.Item(71, 11).Value = 5710
$ws.Cells.Item(71, 12).Value = 14358.7185
$ws.Cells.Item(71, 13).Value = -1966
$ws.Cells.Item(71, 14).Value = -21846.7185

$ws.Cells.Item(138, 8).Value = 152192.6
$ws.Cells.Item(138, 10).Value = 152192.6
$ws.Cells.Item(138, 12).Value = 152192.6
$ws.Cells.Item(138, 14).Value = -162472.6

$ws.Cells.Item(141, 8).Value = 38472.5
$ws.Cells.Item(141, 10).Value = 38472.5
$ws.Cells.Item(141, 12).Value = 38472.5
$ws.Cells.Item(141, 14).Value = -48832.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 353.125
$ws.Cells.Item(107, 9).Value = 275
$ws.Cells.Item(107, 10).Value = 413.8889
$ws.Cells.Item(107, 11).Value = 825
$ws.Cells.Item(107, 12).Value = 1241.6667
$ws.Cells.Item(107, 13).Value = 1095
$ws.Cells.Item(107, 14).Value = -5081.6667

$ws.Cells.Item(136, 8).Value = 733.525
$ws.Cells.Item(136, 9).Value = 679.57574
$ws.Cells.Item(136, 10).Value = 987.8570999999999
$ws.Cells.Item(136, 11).Value = 2038.72722
$ws.Cells.Item(136, 12).Value = 2963.5713
$ws.Cells.Item(136, 13).Value = 511.27278
$ws.Cells.Item(136, 14).Value = -8063.5713
